$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2028169014084507
$ws.Range("C2").Value = 0.5605633802816902
$ws.Range("J2").Value = 0.008450704225352112
$ws.Range("P2").Value = 0.1295774647887324
$ws.Range("S2").Value = 0.09859154929577464
$ws.Range("C3").Value = 0.009900990099009901
$ws.Range("J3").Value = 0.06930693069306931
$ws.Range("P3").Value = 0.6782178217821783
$ws.Range("S3").Value = 0.2425742574257426
$ws.Range("B6").Value = 0.04365079365079365
$ws.Range("D6").Value = 0.0119047619047619
$ws.Range("F6").Value = 0.04761904761904762
$ws.Range("J6").Value = 0.2817460317460317
$ws.Range("O6").Value = 0.03174603174603174
$ws.Range("Q6").Value = 0.123015873015873
$ws.Range("R6").Value = 0.07539682539682539
$ws.Range("S6").Value = 0.3849206349206349
$ws.Range("B7").Value = 0.1421052631578947
$ws.Range("D7").Value = 0.005263157894736842
$ws.Range("F7").Value = 0.06842105263157895
$ws.Range("J7").Value = 0.1684210526315789
$ws.Range("O7").Value = 0.02105263157894737
$ws.Range("Q7").Value = 0.1736842105263158
$ws.Range("R7").Value = 0.06315789473684211
$ws.Range("S7").Value = 0.3578947368421053
$ws.Range("B8").Value = 0.1224899598393574
$ws.Range("D8").Value = 0.01807228915662651
$ws.Range("F8").Value = 0.06224899598393574
$ws.Range("J8").Value = 0.1265060240963855
$ws.Range("O8").Value = 0.02008032128514056
$ws.Range("Q8").Value = 0.1666666666666667
$ws.Range("R8").Value = 0.08433734939759036
$ws.Range("S8").Value = 0.3995983935742972
$ws.Range("B9").Value = 0.1185185185185185
$ws.Range("D9").Value = 0.01481481481481482
$ws.Range("F9").Value = 0.0962962962962963
$ws.Range("J9").Value = 0.08888888888888889
$ws.Range("O9").Value = 0.01481481481481482
$ws.Range("Q9").Value = 0.162962962962963
$ws.Range("R9").Value = 0.1037037037037037
$ws.Range("S9").Value = 0.4
$ws.Range("B10").Value = 0.1202389843166542
$ws.Range("D10").Value = 0.02763256161314414
$ws.Range("E10").Value = 0.003734129947722181
$ws.Range("F10").Value = 0.07916355489171023
$ws.Range("J10").Value = 0.1374159820761763
$ws.Range("O10").Value = 0.0186706497386109
$ws.Range("Q10").Value = 0.1837191934279313
$ws.Range("R10").Value = 0.07692307692307693
$ws.Range("S10").Value = 0.3525018670649739
$ws.Range("G11").Value = 0.1494661921708185
$ws.Range("J11").Value = 0.1067615658362989
$ws.Range("K11").Value = 0.199288256227758
$ws.Range("L11").Value = 0.5266903914590747
$ws.Range("S11").Value = 0.01779359430604982
$ws.Range("G12").Value = 0.7692307692307693
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("K12").Value = 0.00641025641025641
$ws.Range("L12").Value = 0.04487179487179487
$ws.Range("S12").Value = 0.01282051282051282
$ws.Range("G13").Value = 0.6545454545454545
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.07272727272727272
$ws.Range("F15").Value = 0.01746724890829694
$ws.Range("H15").Value = 0.2096069868995633
$ws.Range("I15").Value = 0.08296943231441048
$ws.Range("J15").Value = 0.3318777292576419
$ws.Range("K15").Value = 0.05240174672489083
$ws.Range("M15").Value = 0.008733624454148471
$ws.Range("N15").Value = 0.008733624454148471
$ws.Range("O15").Value = 0.04366812227074236
$ws.Range("S15").Value = 0.2445414847161572
$ws.Range("F16").Value = 0.02816901408450704
$ws.Range("H16").Value = 0.1877934272300469
$ws.Range("I16").Value = 0.06572769953051644
$ws.Range("J16").Value = 0.4178403755868544
$ws.Range("K16").Value = 0.09859154929577464
$ws.Range("M16").Value = 0.02816901408450704
$ws.Range("O16").Value = 0.08450704225352113
$ws.Range("S16").Value = 0.0892018779342723
$ws.Range("F17").Value = 0.01699029126213592
$ws.Range("H17").Value = 0.220873786407767
$ws.Range("I17").Value = 0.05339805825242718
$ws.Range("J17").Value = 0.4029126213592233
$ws.Range("K17").Value = 0.133495145631068
$ws.Range("M17").Value = 0.01941747572815534
$ws.Range("O17").Value = 0.06310679611650485
$ws.Range("S17").Value = 0.08980582524271845
$ws.Range("F18").Value = 0.02162162162162162
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.06486486486486487
$ws.Range("J18").Value = 0.4648648648648649
$ws.Range("K18").Value = 0.0972972972972973
$ws.Range("M18").Value = 0.02702702702702703
$ws.Range("O18").Value = 0.07567567567567568
$ws.Range("S18").Value = 0.04864864864864865
$ws.Range("F19").Value = 0.02245388933440257
$ws.Range("H19").Value = 0.2269446672012831
$ws.Range("I19").Value = 0.05533279871692061
$ws.Range("J19").Value = 0.3889334402566159
$ws.Range("K19").Value = 0.08981555733761026
$ws.Range("M19").Value = 0.02967121090617482
$ws.Range("O19").Value = 0.07217321571772253
$ws.Range("S19").Value = 0.1146752205292702
